$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.514.23'
$ws.Range("E2").Value = '  +1.99%  '
$ws.Range("D3").Value = '1.682.07'
$ws.Range("E3").Value = '  +2.49%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '217.42'
$ws.Range("E5").Value = '  +3.53%  '
$ws.Range("D6").Value = '0.5326'
$ws.Range("E6").Value = '  +2.88%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").Value = '0.2676'
$ws.Range("E8").Value = '  +4.25%  '
$ws.Range("D9").Value = '0.06427'
$ws.Range("E9").Value = '  +2.99%  '
$ws.Range("D10").Value = '21.48'
$ws.Range("E10").Value = '  +5.17%  '
$ws.Range("D11").Value = '0.07793'
$ws.Range("E11").Value = '  +3.29%  '
$ws.Range("D12").Value = '1.696.79'
$ws.Range("E12").Value = '  +3.10%  '
$ws.Range("D13").Value = '4.511'
$ws.Range("E13").Value = '  +3.29%  '
$ws.Range("E14").Value = '  +3.85%  '
$ws.Range("D15").Value = '0.0₅8422'
$ws.Range("E15").Value = '  +5.63%  '
$ws.Range("D16").Value = '65.98'
$ws.Range("E16").Value = '  +1.36%  '
$ws.Range("D17").Value = '26.547.03'
$ws.Range("E17").Value = '  +2.04%  '
$ws.Range("E18").Value = '  -0.01%  '
$ws.Range("D19").Value = '4.805'
$ws.Range("E19").Value = '  +3.40%  '
$ws.Range("D20").Value = '195.51'
$ws.Range("E20").Value = '  +5.44%  '
$ws.Range("E21").Value = '  +3.61%  '
$ws.Range("D22").Value = '6.383'
$ws.Range("E22").Value = '  +4.58%  '
$ws.Range("E23").Value = '  -0.04%  '
$ws.Range("E24").Value = '  -1.46%  '
$ws.Range("D25").Value = '0.1278'
$ws.Range("E25").Value = '  +7.20%  '
$ws.Range("E26").Value = '  +1.54%  '
$ws.Range("E27").Value = '  +4.18%  '
$ws.Range("D28").Value = '1.416'
$ws.Range("E28").Value = '  +3.54%  '
$ws.Range("D29").Value = '0.06127'
$ws.Range("E29").Value = '  +2.43%  '
$ws.Range("E30").Value = '  +2.65%  '
$ws.Range("D31").Value = '3.605'
$ws.Range("E31").Value = '  +7.40%  '
$ws.Range("E32").Value = '  +3.35%  '
$ws.Range("D33").Value = '1.704'
$ws.Range("E33").Value = '  +5.65%  '
$ws.Range("D34").Value = '1.016'
$ws.Range("E34").Value = '  +4.52%  '
$ws.Range("D35").Value = '2.790'
$ws.Range("E35").Value = '  +2.21%  '
$ws.Range("D36").Value = '2.419'
$ws.Range("E36").Value = '  +1.72%  '
$ws.Range("D37").Value = '0.5699'
$ws.Range("E37").Value = '  -2.72%  '
$ws.Range("D38").Value = '0.01641'
$ws.Range("E38").Value = '  +2.87%  '
$ws.Range("D39").Value = '5.954'
$ws.Range("E39").Value = '  +3.68%  '
$ws.Range("D40").Value = '0.8706'
$ws.Range("E40").Value = '  +3.37%  '
$ws.Range("D41").Value = '1.061.45'
$ws.Range("E41").Value = '  +1.40%  '
$ws.Range("E42").Value = '  -0.08%  '
$ws.Range("D43").Value = '99.93'
$ws.Range("E43").Value = '  +0.13%  '
$ws.Range("D44").Value = '1.832.34'
$ws.Range("E44").Value = '  +2.21%  '
$ws.Range("D45").Value = '0.0₈111'
$ws.Range("E45").Value = '  +3.92%  '
$ws.Range("D46").Value = '57.26'
$ws.Range("E46").Value = '  +5.28%  '
$ws.Range("D47").Value = '8.183'
$ws.Range("E47").Value = '  +2.65%  '
$ws.Range("E48").Value = '  -0.48%  '
$ws.Range("D49").Value = '0.05203'
$ws.Range("E49").Value = '  +0.21%  '
$ws.Range("D50").Value = '6.072'
$ws.Range("E50").Value = '  +4.64%  '
$ws.Range("E51").Value = '  +0.28%  '
